$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FACTORS")

# Updated FACTORS_1 (column B) and FACTORS_2 (column C) values for rows 2-12
$values = @{
    2  = @(0.1422929754124884,  0.06206059048872284)
    3  = @(0.07597668412866533, 0.1487447694034979)
    4  = @(0.01824749113468727, 0.07785249618475309)
    5  = @(0.1770615309149006,  0.007590897519115494)
    6  = @(0.1100476951253636,  0.03442464693173394)
    7  = @(0.09189774332174123, 0.1292786773499676)
    8  = @(0.04545117007920132, 0.08896421609904119)
    9  = @(0.1271459160149058,  0.1454769433553602)
    10 = @(0.1511606971817293,  0.02595141623140635)
    11 = @(0.004122276361780318,0.110508344624508)
    12 = @(0.05659582032453687, 0.1691470018118935)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
